# Update the Forecast Portfolio: refresh the 15-minute consumption forecast
# series (columns A & B, rows 2-93) with the new forecast pull, which now
# only spans 91 data points instead of 92 (so the trailing row 93 is removed).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aValues = @(5590, 5540, 5500, 5450, 5410, 5390, 5380, 5370, 5340, 5330, 5340, 5350, 5360, 5380, 5420, 5470, 5520, 5590, 5670, 5760, 5840, 5930, 6020, 6120, 6210, 6320, 6460, 6530, 6550, 6560, 6550, 6540, 6520, 6450, 6390, 6330, 6270, 6220, 6170, 6140, 6110, 6080, 6070, 6060, 6070, 6090, 6110, 6150, 6180, 6230, 6280, 6350, 6440, 6520, 6600, 6700, 6800, 6910, 7040, 7180, 7330, 7460, 7540, 7570, 7560, 7540, 7530, 7500, 7490, 7460, 7420, 7360, 7300, 7240, 7150, 7060, 6970, 6890, 6800, 6700, 6600, 6500, 6420, 6330, 6250, 6180, 6100, 5930, 5880, 5840, 5790)
$bValues = @(46022, 46022.01041666666, 46022.02083333334, 46022.03125, 46022.04166666666, 46022.05208333334, 46022.0625, 46022.07291666666, 46022.08333333334, 46022.09375, 46022.125, 46022.13541666666, 46022.14583333334, 46022.15625, 46022.16666666666, 46022.17708333334, 46022.1875, 46022.19791666666, 46022.20833333334, 46022.21875, 46022.22916666666, 46022.23958333334, 46022.25, 46022.26041666666, 46022.27083333334, 46022.28125, 46022.29166666666, 46022.30208333334, 46022.3125, 46022.32291666666, 46022.34375, 46022.35416666666, 46022.36458333334, 46022.375, 46022.38541666666, 46022.39583333334, 46022.40625, 46022.41666666666, 46022.42708333334, 46022.4375, 46022.44791666666, 46022.45833333334, 46022.46875, 46022.47916666666, 46022.5, 46022.51041666666, 46022.52083333334, 46022.53125, 46022.54166666666, 46022.55208333334, 46022.5625, 46022.57291666666, 46022.58333333334, 46022.59375, 46022.60416666666, 46022.61458333334, 46022.625, 46022.63541666666, 46022.64583333334, 46022.65625, 46022.66666666666, 46022.67708333334, 46022.6875, 46022.69791666666, 46022.71875, 46022.72916666666, 46022.73958333334, 46022.75, 46022.76041666666, 46022.77083333334, 46022.78125, 46022.79166666666, 46022.80208333334, 46022.8125, 46022.82291666666, 46022.83333333334, 46022.84375, 46022.85416666666, 46022.86458333334, 46022.875, 46022.88541666666, 46022.89583333334, 46022.90625, 46022.91666666666, 46022.92708333334, 46022.9375, 46022.94791666666, 46022.95833333334, 46022.96875, 46022.97916666666, 46022.98958333334)

for ($i = 0; $i -lt $aValues.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $aValues[$i]
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}

$ws.Rows.Item(93).Delete()
